$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.105.04'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '3.158.70'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('D5').Value = '''601.57'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').Value = '''153.66'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').Value = '''0.550'
$ws.Range('E8').Value = '  +3.27%  '
$ws.Range('D9').Value = '3.154.57'
$ws.Range('E9').Value = '  -1.34%  '
$ws.Range('E10').Value = '  -1.92%  '
$ws.Range('D11').Value = '''5.50'
$ws.Range('E11').Value = '  -10.62%  '
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').Value = '''0.0000265'
$ws.Range('E13').Value = '  -2.16%  '
$ws.Range('D14').Value = '''38.26'
$ws.Range('E14').Value = '  -0.51%  '
$ws.Range('D15').Value = '3.679.64'
$ws.Range('E15').Value = '  -1.30%  '
$ws.Range('D16').Value = '66.210.46'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').Value = '3.160.88'
$ws.Range('E18').Value = '  -1.30%  '
$ws.Range('E19').Value = '  +1.12%  '
$ws.Range('D20').Value = '''508.04'
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D21').Value = '''15.36'
$ws.Range('E21').Value = '  -1.30%  '
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('E23').Value = '  +1.06%  '
$ws.Range('E24').Value = '  -3.81%  '
$ws.Range('D25').Value = '''84.44'
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  -0.77%  '
$ws.Range('E28').Value = '  -1.14%  '
$ws.Range('D29').Value = '''2.38'
$ws.Range('E29').Value = '  +6.45%  '
$ws.Range('D30').Value = '''3.03'
$ws.Range('E30').Value = '  +6.00%  '
$ws.Range('D31').Value = '''6.90'
$ws.Range('E31').Value = '  +0.63%  '
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').Value = '''1.19'
$ws.Range('E34').Value = '  -1.59%  '
$ws.Range('D35').Value = '''6.47'
$ws.Range('E35').Value = '  -1.93%  '
$ws.Range('D36').Value = '''503.49'
$ws.Range('E36').Value = '  +4.51%  '
$ws.Range('D37').Value = '''54.73'
$ws.Range('E37').Value = '  -0.95%  '
$ws.Range('D38').Value = '''0.0881'
$ws.Range('E38').Value = '  -2.86%  '
$ws.Range('E39').Value = '  -0.60%  '
$ws.Range('D40').Value = '''0.128'
$ws.Range('E40').Value = '  +7.68%  '
$ws.Range('D42').Value = '0.0₃0673'
$ws.Range('E42').Value = '  +5.79%  '
$ws.Range('D43').Value = '''0.294'
$ws.Range('E43').Value = '  -1.12%  '
$ws.Range('D44').Value = '''2.77'
$ws.Range('E44').Value = '  -6.35%  '
$ws.Range('D45').Value = '''2.40'
$ws.Range('E45').Value = '  -2.40%  '
$ws.Range('D46').Value = '2.818.55'
$ws.Range('E46').Value = '  -4.29%  '
$ws.Range('D47').Value = '''27.75'
$ws.Range('E47').Value = '  -3.20%  '
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('E49').Value = '  +1.52%  '
$ws.Range('E50').Value = '  +0.58%  '
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').Value = '''2.57'
$ws.Range('E51').Value = '  +5.54%  '
